$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.247.98"
$ws.Range("E2").Value = "  +5.45%  "

$ws.Range("D3").Value = "3.454.49"
$ws.Range("E3").Value = "  +6.28%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.58"
$ws.Range("E5").Value = "  +6.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.53"
$ws.Range("E6").Value = "  +7.29%  "

$ws.Range("D8").Value = "3.460.12"
$ws.Range("E8").Value = "  +6.07%  "

$ws.Range("E9").Value = "  +1.31%  "

$ws.Range("E10").Value = "  +2.96%  "

$ws.Range("E11").Value = "  +7.04%  "

$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("D13").Value = "4.049.14"
$ws.Range("E13").Value = "  +6.19%  "

$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("E15").Value = "  +7.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.62"
$ws.Range("E16").Value = "  +5.12%  "

$ws.Range("D17").Value = "64.264.11"
$ws.Range("E17").Value = "  +5.50%  "

$ws.Range("D18").Value = "3.452.03"
$ws.Range("E18").Value = "  +6.10%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  +2.36%  "

$ws.Range("E20").Value = "  +7.39%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "397.27"
$ws.Range("E21").Value = "  +5.34%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.56"
$ws.Range("E22").Value = "  +1.92%  "

$ws.Range("E23").Value = "  +2.61%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.14"
$ws.Range("E25").Value = "  +3.05%  "

$ws.Range("E26").Value = "  +19.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.50"
$ws.Range("E27").Value = "  +10.26%  "

$ws.Range("E28").Value = "  +6.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("E30").Value = "  +13.64%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("E31").Value = "  +8.79%  "

$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.70"
$ws.Range("E32").Value = "  +8.19%  "

$ws.Range("E33").Value = "  +6.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.79"
$ws.Range("E34").Value = "  +5.26%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.87"
$ws.Range("E36").Value = "  +3.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  +4.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.62"
$ws.Range("E38").Value = "  -0.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.30"
$ws.Range("E39").Value = "  +7.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0786"
$ws.Range("E40").Value = "  +9.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  +9.16%  "

$ws.Range("D42").Value = "2.861.69"
$ws.Range("E42").Value = "  +1.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0322"
$ws.Range("E43").Value = "  +2.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.781"
$ws.Range("E44").Value = "  +6.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.18"
$ws.Range("E45").Value = "  +5.33%  "

$ws.Range("E46").Value = "  +2.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  +9.95%  "

$ws.Range("D48").Value = "3.501.00"
$ws.Range("E48").Value = "  +6.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.74"
$ws.Range("E49").Value = "  +6.11%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "298.94"
$ws.Range("E50").Value = "  +8.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.12"
$ws.Range("E51").Value = "  +22.44%  "
